# Append the 2025-10-20 profit row (row 64) to the bottom of the profit log.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds dates formatted/stored as plain text (matching the rest of
# the sheet's "MM/DD/YYYY" text cells), not an Excel date serial. Force the
# cell to text via NumberFormat before assigning so "10/20/2025" isn't
# auto-converted to a date value, then clear the forced format so the cell
# keeps the sheet's default (unstyled) appearance like its neighbours.
$rowNum = 64
$ws.Range("A$rowNum").NumberFormat = "@"
$ws.Range("A$rowNum").Value = "10/20/2025"
$ws.Range("A$rowNum").ClearFormats()

# Column B is the numeric profit figure for that date.
$ws.Range("B$rowNum").Value = 10101.96
